# Weekly CompStat (115th Precinct) data refresh:
#  - Bump "Volume ... Number" from 46 to 47
#  - Move the reporting week forward one week (11/13-11/19/2023 -> 11/20-11/26/2023)
#  - Refresh the Crime Complaints table (rows 14-30) with the newly collected figures,
#    including a handful of cells that flip between a numeric count and the "0" / "***.*"
#    placeholder text used when a rate cannot be computed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (these live inside multi-run shared strings; editing the
# relevant substring via Characters() keeps the rest of the sentence intact).
# ---------------------------------------------------------------------------

$volumeCell = $ws.Range("A8")
$volumeText = $volumeCell.Value()
$numIdx = $volumeText.LastIndexOf("46")
$volumeCell.Characters($numIdx + 1, 2).Text() = "47"

$weekCell = $ws.Range("C9")

$weekText = $weekCell.Value()
$fromIdx = $weekText.IndexOf("11/13/2023")
$weekCell.Characters($fromIdx + 1, 10).Text() = "11/20/2023"

$weekText = $weekCell.Value()
$throughIdx = $weekText.IndexOf("11/19/2023")
$weekCell.Characters($throughIdx + 1, 10).Text() = "11/26/2023"

# ---------------------------------------------------------------------------
# Donor cells used purely as a source of "already correct" styling, so that
# cells which switch between numeric and placeholder-text ("0" / "***.*")
# pick up the exact formatting Excel would apply in each case. None of these
# donor cells are themselves modified by this script.
# ---------------------------------------------------------------------------

$donorText0 = $ws.Range("D14")     # General-formatted text cell currently showing "0"
$donorTextStar = $ws.Range("E14")  # General-formatted text cell currently showing "***.*"
$donorNum15 = $ws.Range("F14")     # #,##0-formatted numeric cell
$donorNum16 = $ws.Range("K14")     # percent-style numeric cell

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
$donorText0.Copy($ws.Range("C14"))
$ws.Range("N14").Value() = -82.352941176470

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$donorText0.Copy($ws.Range("C15"))
$donorText0.Copy($ws.Range("D15"))
$donorTextStar.Copy($ws.Range("E15"))
$ws.Range("G15").Value() = 1
$ws.Range("H15").Value() = 0
$ws.Range("L15").Value() = -25.806451612903
$ws.Range("N15").Value() = -28.125

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value() = 7
$ws.Range("D16").Value() = 9
$ws.Range("E16").Value() = -22.222222222222
$ws.Range("F16").Value() = 29
$ws.Range("G16").Value() = 37
$ws.Range("H16").Value() = -21.621621621621
$ws.Range("I16").Value() = 324
$ws.Range("J16").Value() = 286
$ws.Range("K16").Value() = 13.286713286713
$ws.Range("L16").Value() = 47.272727272727
$ws.Range("M16").Value() = 0.934579439252
$ws.Range("N16").Value() = -72.795969773299

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value() = 10
$ws.Range("D17").Value() = 5
$ws.Range("E17").Value() = 100
$ws.Range("G17").Value() = 24
$ws.Range("H17").Value() = 33.333333333333
$ws.Range("I17").Value() = 435
$ws.Range("J17").Value() = 379
$ws.Range("K17").Value() = 14.775725593667
$ws.Range("L17").Value() = 24.285714285714
$ws.Range("M17").Value() = 45.973154362416
$ws.Range("N17").Value() = 11.253196930946

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value() = 1
$ws.Range("E18").Value() = -50
$ws.Range("G18").Value() = 10
$ws.Range("H18").Value() = -20
$ws.Range("I18").Value() = 129
$ws.Range("J18").Value() = 135
$ws.Range("K18").Value() = -4.444444444444
$ws.Range("L18").Value() = 2.380952380952
$ws.Range("M18").Value() = -52.044609665427
$ws.Range("N18").Value() = -92.825361512792

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value() = 15
$ws.Range("D19").Value() = 20
$ws.Range("E19").Value() = -25
$ws.Range("F19").Value() = 63
$ws.Range("G19").Value() = 76
$ws.Range("H19").Value() = -17.105263157894
$ws.Range("I19").Value() = 771
$ws.Range("J19").Value() = 891
$ws.Range("K19").Value() = -13.468013468013
$ws.Range("L19").Value() = 35.978835978836
$ws.Range("M19").Value() = 69.078947368421
$ws.Range("N19").Value() = -41.457858769931

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value() = 6
$ws.Range("D20").Value() = 7
$ws.Range("E20").Value() = -14.285714285714
$ws.Range("I20").Value() = 309
$ws.Range("J20").Value() = 288
$ws.Range("K20").Value() = 7.291666666666
$ws.Range("L20").Value() = 63.492063492063
$ws.Range("M20").Value() = 45.070422535211
$ws.Range("N20").Value() = -85.021812893843

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value() = 39
$ws.Range("D21").Value() = 43
$ws.Range("E21").Value() = -9.302325581395
$ws.Range("F21").Value() = 159
$ws.Range("G21").Value() = 176
$ws.Range("H21").Value() = -9.659090909090
$ws.Range("I21").Value() = 1994
$ws.Range("J21").Value() = 2018
$ws.Range("K21").Value() = -1.189296333002
$ws.Range("L21").Value() = 34.095494283792
$ws.Range("M21").Value() = 25.645872715816
$ws.Range("N21").Value() = -70.715229842855

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$donorNum15.Copy($ws.Range("C22"))
$ws.Range("C22").Value() = 2
$ws.Range("G22").Value() = 4
$ws.Range("H22").Value() = 0
$ws.Range("I22").Value() = 78
$ws.Range("K22").Value() = 34.482758620689
$ws.Range("L22").Value() = 188.888888888889
$ws.Range("M22").Value() = 178.571428571429

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value() = 34
$ws.Range("D24").Value() = 40
$ws.Range("E24").Value() = -15
$ws.Range("F24").Value() = 157
$ws.Range("G24").Value() = 203
$ws.Range("H24").Value() = -22.660098522167
$ws.Range("I24").Value() = 1799
$ws.Range("J24").Value() = 1770
$ws.Range("K24").Value() = 1.638418079096
$ws.Range("L24").Value() = 43.804956035171
$ws.Range("M24").Value() = 74.660194174757

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value() = 24
$ws.Range("D25").Value() = 19
$ws.Range("E25").Value() = 26.315789473684
$ws.Range("F25").Value() = 88
$ws.Range("H25").Value() = 23.943661971831
$ws.Range("I25").Value() = 895
$ws.Range("J25").Value() = 822
$ws.Range("K25").Value() = 8.880778588807
$ws.Range("L25").Value() = 15.038560411311
$ws.Range("M25").Value() = 5.791962174940

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
$donorText0.Copy($ws.Range("C26"))
$donorText0.Copy($ws.Range("D26"))
$donorTextStar.Copy($ws.Range("E26"))
$ws.Range("G26").Value() = 2
$ws.Range("H26").Value() = 50
$ws.Range("L26").Value() = -11.764705882352

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C27").Value() = 1
$donorNum15.Copy($ws.Range("D27"))
$ws.Range("D27").Value() = 2
$donorNum16.Copy($ws.Range("E27"))
$ws.Range("E27").Value() = -50
$ws.Range("F27").Value() = 7
$ws.Range("G27").Value() = 18
$ws.Range("H27").Value() = -61.111111111111
$ws.Range("I27").Value() = 133
$ws.Range("J27").Value() = 104
$ws.Range("K27").Value() = 27.884615384615
$ws.Range("L27").Value() = 38.541666666666

# ---------------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------------
$donorNum15.Copy($ws.Range("C30"))
$ws.Range("C30").Value() = 1
